$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.613.48'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '1.912.01'
$ws.Range("E3").Value = '  +3.90%  '
$ws.Range("E4").Value = '  +0.59%  '
$ws.Range("D5").Value = "'245.14"
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").Value = "'42.42"
$ws.Range("E8").Value = '  +3.50%  '
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("D10").Value = "'0.0708"
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("D11").Value = "'0.0996"
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").Value = '2.188.28'
$ws.Range("E12").Value = '  +3.87%  '
$ws.Range("D13").Value = "'12.54"
$ws.Range("E13").Value = '  +10.58%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = "'0.694"
$ws.Range("E14").Value = '  +3.64%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = "'4.85"
$ws.Range("E15").Value = '  +4.33%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.872.52'
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").Value = '35.604.42'
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("D18").Value = "'72.07"
$ws.Range("E18").Value = '  +3.23%  '
$ws.Range("D19").Value = '0.0₃0810'
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").Value = "'243.95"
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = "'12.52"
$ws.Range("E21").Value = '  +3.29%  '
$ws.Range("E22").Value = '  +4.03%  '
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").Value = "'2.29"
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").Value = "'171.44"
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").Value = "'2.11"
$ws.Range("E26").Value = '  +28.25%  '
$ws.Range("E27").Value = '  +8.17%  '
$ws.Range("D28").Value = "'18.03"
$ws.Range("E28").Value = '  +3.61%  '
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("D30").Value = "'4.10"
$ws.Range("E30").Value = '  +4.02%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'0.952"
$ws.Range("E31").Value = '  +27.48%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = "'0.0566"
$ws.Range("E32").Value = '  +2.90%  '
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("E34").Value = '  +6.54%  '
$ws.Range("E35").Value = '  +6.97%  '
$ws.Range("E36").Value = '  +5.22%  '
$ws.Range("E37").Value = '  +5.69%  '
$ws.Range("E38").Value = '  +5.39%  '
$ws.Range("E39").Value = '  +4.45%  '
$ws.Range("D40").Value = "'91.89"
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("D41").Value = '1.362.73'
$ws.Range("E41").Value = '  +1.86%  '
$ws.Range("D42").Value = "'15.21"
$ws.Range("E42").Value = '  +4.45%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = "'49.16"
$ws.Range("E43").Value = '  +44.39%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = "'0.0593"
$ws.Range("E44").Value = '  +11.70%  '
$ws.Range("D45").Value = "'13.06"
$ws.Range("E45").Value = '  +19.68%  '
$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = '  +4.49%  '
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("E48").Value = '  +1.27%  '
$ws.Range("E49").Value = '  +5.89%  '
$ws.Range("D50").Value = '2.096.84'
$ws.Range("E50").Value = '  +3.46%  '
$ws.Range("E51").Value = '  +5.14%  '
